$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row that held the outlier "v" measurement (old row 25, element 54)
$ws.Rows.Item(25).Delete()
# Delete the trailing "p" measurement row (old row 18, element 65)
$ws.Rows.Item(18).Delete()

# Refresh the "p" (power) measurement values with the newly re-run simulation results
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = -0.035
$ws.Range("E4").Value = -0.03
$ws.Range("E5").Value = -0.04
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = -0.045
$ws.Range("E8").Value = -0.065
$ws.Range("E9").Value = -0.015
$ws.Range("E10").Value = -0.05
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = -0.01225
$ws.Range("E13").Value = -0.0175
$ws.Range("E14").Value = -0.015399999999999999
$ws.Range("E15").Value = -0.01225
$ws.Range("E16").Value = -0.01435
$ws.Range("E17").Value = 0

# Renumber the "v" (voltage) measurement element indices sequentially (48-64)
$ws.Range("D18").Value = 48
$ws.Range("D19").Value = 49
$ws.Range("D20").Value = 50
$ws.Range("D21").Value = 51
$ws.Range("D22").Value = 52
$ws.Range("D23").Value = 53
$ws.Range("D24").Value = 54
$ws.Range("D25").Value = 55
$ws.Range("D26").Value = 56
$ws.Range("D27").Value = 57
$ws.Range("D28").Value = 58
$ws.Range("D29").Value = 59
$ws.Range("D30").Value = 60
$ws.Range("D31").Value = 61
$ws.Range("D32").Value = 62
$ws.Range("D33").Value = 63
$ws.Range("D34").Value = 64

# Restore the active selection to match the saved workbook state
$ws.Range("K10").Select()
